# Add a new worksheet "24AA025UID" (register map for the 24AA025UID EEPROM's
# factory-programmed ID page) right after "IOExpander" and before "SPI".
$wb = $excel.ActiveWorkbook

$newSheet = $wb.Worksheets.Add($null, $wb.Worksheets.Item("IOExpander"))
$newSheet.Name = "24AA025UID"

# --- Header row (reuses the existing shared strings from the other sheets) ---
$newSheet.Range("A1").Value = "Name"
$newSheet.Range("B1").Value = "Hex Address"
$newSheet.Range("C1").Value = "Default Value"
$newSheet.Range("D1").Value = "Bit Width"
$newSheet.Range("E1").Value = "Bit Index (High)"
$newSheet.Range("F1").Value = "Bit Index (Low)"

# --- Column A: register names ---
$newSheet.Range("A2").Value = "SERIAL_NUMBER"
$newSheet.Range("A3").Value = "MANUFACTURER_CODE"
$newSheet.Range("A4").Value = "DEVICE_CODE"

# --- Column C: default values ---
$newSheet.Range("C2").Value = "0x00000000"
$newSheet.Range("C3").Value = "0x29"
$newSheet.Range("C4").Value = "0x41"

# --- Column B: hex addresses ---
$newSheet.Range("B2").Value = "0xfc"
$newSheet.Range("B3").Value = "0xfa"
$newSheet.Range("B4").Value = "0xfb"

# --- Column D: bit widths ---
$newSheet.Range("D2").Value = 32
$newSheet.Range("D3").Value = 8
$newSheet.Range("D4").Value = 8

# --- Column E: bit index (high) ---
$newSheet.Range("E2").Value = "None"
$newSheet.Range("E3").Value = "None"
$newSheet.Range("E4").Value = "None"

# --- Column F: bit index (low) ---
$newSheet.Range("F2").Value = "None"
$newSheet.Range("F3").Value = "None"
$newSheet.Range("F4").Value = "None"

# Highlight the 32-bit default value with the built-in "Neutral" cell style.
$newSheet.Range("C2").Style = "Neutral"

# Column widths to roughly match the source sheet's auto-fit content widths.
$newSheet.Columns.Item(1).ColumnWidth = 20.89
$newSheet.Columns.Item(2).ColumnWidth = 13.44
$newSheet.Columns.Item(3).ColumnWidth = 16.07
$newSheet.Columns.Item(5).ColumnWidth = 15.17
$newSheet.Columns.Item(6).ColumnWidth = 14.44

# Make the new sheet the active tab, with B3 selected.
$newSheet.Activate() | Out-Null
$newSheet.Range("B3").Select() | Out-Null

# The header row on IOExpander was selected (copied) while building this sheet.
$ioExpander = $wb.Worksheets.Item("IOExpander")
$ioExpander.Activate() | Out-Null
$ioExpander.Range("A1:F1").Select() | Out-Null

# Leave the new sheet as the active tab.
$newSheet.Activate() | Out-Null
